$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump the expected qlVersion() result to 1.8.2 so the unit test (B3 vs C3) passes.
$ws.Range("B3").Value = "1.8.2"
